$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.545
$ws.Range("B3").Value = 5.688000000000001
$ws.Range("E3").Value = 16.261
$ws.Range("E12").Value = 17.724
$ws.Range("B14").Value = 5.735
$ws.Range("B16").Value = 5.449
$ws.Range("C18").Value = -11.87
$ws.Range("B21").Value = 9.49
$ws.Range("B23").Value = 7.725
$ws.Range("C24").Value = -12.74
$ws.Range("E24").Value = 16.987
$ws.Range("B25").Value = 5.403
$ws.Range("C25").Value = -11.52
$ws.Range("E25").Value = 17.082
$ws.Range("B26").Value = 6.461
$ws.Range("C27").Value = -13.645
$ws.Range("B29").Value = 5.383999999999999
$ws.Range("C30").Value = -13.165
$ws.Range("C31").Value = -13.223
$ws.Range("C39").Value = -12.728
$ws.Range("B40").Value = 9.269
$ws.Range("E41").Value = 16.473
$ws.Range("C42").Value = -12.67
$ws.Range("C48").Value = -11.52
$ws.Range("E50").Value = 16.305
$ws.Range("C51").Value = -11.02
$ws.Range("C52").Value = -11.684
$ws.Range("B53").Value = 6.687
$ws.Range("E53").Value = 17.041
$ws.Range("C55").Value = -13.618
$ws.Range("C56").Value = -13.445
$ws.Range("E56").Value = 16.02
$ws.Range("B57").Value = 4.864000000000001
$ws.Range("C57").Value = -13.852
$ws.Range("E57").Value = 16.423
$ws.Range("E58").Value = 16.545
$ws.Range("B59").Value = 4.714
$ws.Range("C60").Value = -13.048
$ws.Range("E61").Value = 16.721
$ws.Range("E63").Value = 17.456
$ws.Range("E64").Value = 17.481
$ws.Range("B65").Value = 5.936000000000001
$ws.Range("B69").Value = 5.384
$ws.Range("E70").Value = 17.602
$ws.Range("E72").Value = 17.143
$ws.Range("C73").Value = -12.752
$ws.Range("C74").Value = -12.453
$ws.Range("B79").Value = 5.6
$ws.Range("B83").Value = 5.45
$ws.Range("E86").Value = 16.249
$ws.Range("C89").Value = -10.693
$ws.Range("E89").Value = 17.187
$ws.Range("C90").Value = -13.298
$ws.Range("B91").Value = 5.540999999999999
$ws.Range("C92").Value = -10.927
$ws.Range("B93").Value = 5.472999999999999
$ws.Range("E98").Value = 16.183
$ws.Range("B100").Value = 5.051
$ws.Range("E100").Value = 16.537
$ws.Range("E102").Value = 16.318
